$wb = $excel.ActiveWorkbook

# --- Update data values on the "ranges" sheet ---
$ranges = $wb.Worksheets.Item("ranges")
$ranges.Range("F5").Value = 1
$ranges.Range("E10").Value = 0.5
$ranges.Range("F10").Value = 1.5
$ranges.Range("D11").Value = 0.005
$ranges.Range("F11").Value = 0.01

# --- Update the selection on the "incubation" sheet (no longer the active tab) ---
$incubation = $wb.Worksheets.Item("incubation")
$incubation.Activate() | Out-Null
$incubation.Range("A12:H12").Select() | Out-Null

# --- Update the selection on the "field" sheet ---
$field = $wb.Worksheets.Item("field")
$field.Activate() | Out-Null
$field.Range("A12:H12").Select() | Out-Null

# --- Update the selection on the "ranges" sheet and make it the active tab ---
$ranges.Activate() | Out-Null
$ranges.Range("E11").Select() | Out-Null

Write-Output "applied edits"
